$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F header is "dSF". Update specific cells per repull of data.
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F11").Value = -2
